# Atualização da validação do import excel com os novos campos
# Adds two new columns (tipo_produto, outros_detalhes) to the product import template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("P1").Value = "tipo_produto"
$ws.Range("Q1").Value = "outros_detalhes"

# Data rows - tipo_produto
$ws.Range("P2").Value = "Insumos"
$ws.Range("P3").Value = "Imobilizado"
$ws.Range("P4").Value = "Revenda"
$ws.Range("P5").Value = "Revenda"
$ws.Range("P6").Value = "Revenda"
$ws.Range("P7").Value = "Revenda"

# Data rows - outros_detalhes
$ws.Range("Q4").Value = "sem detalhes"

# Copy the row style (red font) used by the rest of the data rows onto the new cells
$ws.Range("P2:P7").Font.Color = $ws.Range("A2").Font.Color
$ws.Range("Q4").Font.Color = $ws.Range("A2").Font.Color

# Update the view to match the saved state (scrolled right, new cell selected)
$ws.Range("Q5").Select()
$excel.ActiveWindow.ScrollColumn = 4
